# Add season record columns (Wins / Losses / Ties) to the roster sheet.
#
# The sheet currently holds team/player stats in columns A:AC (rows 1-53,
# where row 1 is the header). We append three new columns: AD (Wins),
# AE (Losses), AF (Ties) with the same header style as the rest of row 1,
# and fill every data row (2-53) with the season record (92-70-0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 53

# --- Header row ---------------------------------------------------------
# Copy the formatting of the last existing header cell (AC1, bold + border)
# onto the three new header cells so they match the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# --- Data rows -----------------------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value2 = 92   # AD - Wins
    $ws.Cells.Item($r, 31).Value2 = 70   # AE - Losses
    $ws.Cells.Item($r, 32).Value2 = 0    # AF - Ties
}
